{"js": "// Apply the diff to the single table in the document.\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount\");\nawait context.sync();\n\n// --- Simple in-place text replacements (row index -> new text) ---\n// Row 0: \"99.85\" -> \"0M\"\ntable.getCell(0, 0).value = \"0M\";\n// Row 1: \"0.22\" -> \"0M\"\ntable.getCell(1, 0).value = \"0M\";\n// Row 2: \"152\" -> \"0M\"\ntable.getCell(2, 0).value = \"0M\";\n// Row 3: \"869\" -> \"1869\"\ntable.getCell(3, 0).value = \"1869\";\nawait context.sync();\n\n// --- Insert 3 new rows right after row 3 (index 3), matching the diff's\n//     added <w:tr> blocks (0.00002 / 0.00054 / 0.00010) ---\nconst rowsForInsert = table.rows;\nrowsForInsert.load(\"items\");\nawait context.sync();\nrowsForInsert.items[3].insertRows(\"After\", 3, [[\"0.00002\"], [\"0.00054\"], [\"0.00010\"]]);\nawait context.sync();\n\n// After the insert, the table rows shift:\n//   old row4 \"0.00003\" is now at index 7\n//   old row5 \"0.00042\" is now at index 8\n//   old row6 \"0.00008\" is now at index 9\n//   old row7 \"0.00004\" is now at index 10\n//   old row8 \"0.00009\" is now at index 11\n//   old row9 \"0.00010\" is now at index 12\n//   old row10 \"0.00010\" is now at index 13\n//   old row11 \"0.07378\" is now at index 14\ntable.getCell(8, 0).value = \"0.00018\";   // was \"0.00042\"\ntable.getCell(9, 0).value = \"0.00019\";   // was \"0.00008\"\ntable.getCell(10, 0).value = \"0.00021\"; // was \"0.00004\"\ntable.getCell(11, 0).value = \"0.22400\"; // was \"0.00009\"\nawait context.sync();\n\n// Delete the 3 rows that followed (old 0.00010 / 0.00010 / 0.07378), which\n// are now sitting right after the row we just rewrote to \"0.22400\"\n// (index 12, deleted one at a time with a sync in between so each delete\n// targets the freshly-shifted row 12).\nlet delRows = table.rows;\ndelRows.load(\"items\");\nawait context.sync();\ndelRows.items[12].delete();\nawait context.sync();\n\ndelRows = table.rows;\ndelRows.load(\"items\");\nawait context.sync();\ndelRows.items[12].delete();\nawait context.sync();\n\ndelRows = table.rows;\ndelRows.load(\"items\");\nawait context.sync();\ndelRows.items[12].delete();\nawait context.sync();\n\n// --- Collapse the final 3 multi-run rows down to single values ---\ntable.load(\"rowCount\");\nawait context.sync();\nconst n = table.rowCount;\ntable.getCell(n - 3, 0).value = \"99.85\";\ntable.getCell(n - 2, 0).value = \"0.22\";\ntable.getCell(n - 1, 0).value = \"152\";\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# --- Simple in-place text replacements (1-based row index) ---\n$t.Cell(1,1).Range.Text = \"0M\"      # was \"99.85\"\n$t.Cell(2,1).Range.Text = \"0M\"      # was \"0.22\"\n$t.Cell(3,1).Range.Text = \"0M\"      # was \"152\"\n$t.Cell(4,1).Range.Text = \"1869\"    # was \"869\"\n\n# --- Insert 3 new rows right before row 5 (the old \"0.00003\" row),\n#     matching the diff's added <w:tr> blocks. Each Rows.Add(beforeRow) call\n#     inserts immediately above beforeRow, so insert in REVERSE order\n#     (last desired row first) to end up with 0.00002 / 0.00054 / 0.00010\n#     reading top-to-bottom. ---\n$beforeRow = $t.Rows.Item(5)\n$newRow3 = $t.Rows.Add($beforeRow)\n$newRow3.Cells.Item(1).Range.Text = \"0.00010\"\n$newRow2 = $t.Rows.Add($beforeRow)\n$newRow2.Cells.Item(1).Range.Text = \"0.00054\"\n$newRow1 = $t.Rows.Add($beforeRow)\n$newRow1.Cells.Item(1).Range.Text = \"0.00002\"\n\n# After the 3 inserts, the old rows shift down by 3:\n#   old row5  \"0.00003\" -> row 8  (unchanged)\n#   old row6  \"0.00042\" -> row 9\n#   old row7  \"0.00008\" -> row 10\n#   old row8  \"0.00004\" -> row 11\n#   old row9  \"0.00009\" -> row 12\n#   old row10 \"0.00010\" -> row 13\n#   old row11 \"0.00010\" -> row 14\n#   old row12 \"0.07378\" -> row 15\n$t.Cell(9,1).Range.Text  = \"0.00018\"  # was \"0.00042\"\n$t.Cell(10,1).Range.Text = \"0.00019\"  # was \"0.00008\"\n$t.Cell(11,1).Range.Text = \"0.00021\"  # was \"0.00004\"\n$t.Cell(12,1).Range.Text = \"0.22400\"  # was \"0.00009\"\n\n# Delete the 3 rows that used to hold 0.00010 / 0.00010 / 0.07378 \u2014 they now\n# sit immediately after the row we just rewrote (row 13, three times).\n$t.Rows.Item(13).Delete()\n$t.Rows.Item(13).Delete()\n$t.Rows.Item(13).Delete()\n\n# --- Collapse the final 3 multi-run rows down to single values ---\n$n = $t.Rows.Count\n$t.Cell($n-2,1).Range.Text = \"99.85\"\n$t.Cell($n-1,1).Range.Text = \"0.22\"\n$t.Cell($n,1).Range.Text   = \"152\"\n"}
